$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I6").Value = "jdavid@edeq.com.co"
$ws.Range("J6").Value = "lamariad@gmail.com"
$ws.Range("I5").Value = "jdavid"
$ws.Range("J5").Value = "lamariad"
